$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the prompt labels from "cost" to "price"
$ws.Range("D6").Value = "Enter price of 1 pizza:"
$ws.Range("H6").Value = "Enter price of 1 book:"

# Update the price of 1 pizza from 12 to 8
$ws.Range("F6").Value = 8

# Update the active cell selection to match the recorded state
$ws.Range("H7").Select()
